$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Basic Movement
$ws.Range("B10").Value = "Basic Movement"
$ws.Range("C10").Value = "Xavier Trillo"
$ws.Range("D10").Value = "1 hour and 30 minutes"
$ws.Range("E10").Value = "2 hours"

# Row 11: Smoth Animations
$ws.Range("B11").Value = "Smoth Animations"
$ws.Range("C11").Value = "Xavier Trillo"
$ws.Range("D11").Value = "0 hours and 30 minutes"
$ws.Range("E11").Value = "1 hour"

# Row 12: Collision Boxes
$ws.Range("B12").Value = "Collision Boxes"
$ws.Range("C12").Value = "Xavier Trillo"
$ws.Range("D12").Value = "0 hours and 30 minutes"
$ws.Range("E12").Value = "0 hours and 30 minutes"

$ws.Range("E12").Select()
